$d = $word.ActiveDocument

$pairs = @(
    @("96×61=", "97×46="),
    @("29×21=", "39×77="),
    @("57×82=", "76×59="),
    @("38×18=", "45×41="),
    @("48×48=", "45×38="),
    @("90×41=", "11×46="),
    @("48×93=", "82×91="),
    @("97×73=", "93×92="),
    @("83×63=", "73×31="),
    @("72×93=", "20×17="),
    @("36×18=", "32×95="),
    @("26×98=", "24×77="),
    @("37×96=", "60×78="),
    @("19×46=", "70×56="),
    @("42×90=", "36×33="),
    @("25×71=", "33×92="),
    @("80×23=", "79×52="),
    @("32×88=", "62×24="),
    @("67×25=", "73×62="),
    @("36×56=", "47×71="),
    @("80×79=", "18×30="),
    @("50×27=", "38×79="),
    @("74×85=", "96×23="),
    @("79×51=", "87×32="),
    @("21×45=", "93×53=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
